$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.106.71'
$ws.Range('E2').Value = '  -1.13%  '
$ws.Range('D3').Value = '2.630.04'
$ws.Range('E3').Value = '  -2.11%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '594.20'
$ws.Range('E5').Value = '  -1.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.65'
$ws.Range('E6').Value = '  +0.97%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.534'
$ws.Range('E8').Value = '  -2.74%  '
$ws.Range('D9').Value = '2.629.61'
$ws.Range('E9').Value = '  -1.99%  '
$ws.Range('E10').Value = '  -1.74%  '
$ws.Range('E11').Value = '  +1.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.359'
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.23'
$ws.Range('E13').Value = '  -0.19%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.75'
$ws.Range('E14').Value = '  -0.74%  '
$ws.Range('D15').Value = '3.111.41'
$ws.Range('E15').Value = '  -1.95%  '
$ws.Range('E16').Value = '  -2.10%  '
$ws.Range('D17').Value = '67.040.21'
$ws.Range('E17').Value = '  -1.12%  '
$ws.Range('D18').Value = '2.629.66'
$ws.Range('E18').Value = '  -2.73%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.18'
$ws.Range('E19').Value = '  +4.43%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.07'
$ws.Range('E20').Value = '  +7.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '357.35'
$ws.Range('E21').Value = '  -2.55%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.32'
$ws.Range('E22').Value = '  -2.52%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.67'
$ws.Range('E23').Value = '  -3.59%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.91'
$ws.Range('E24').Value = '  +10.09%  '
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('E26').Value = '  -5.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '70.24'
$ws.Range('D28').Value = '2.757.67'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  -0.22%  '
$ws.Range('E30').Value = '  -1.72%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '550.39'
$ws.Range('E31').Value = '  -2.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.92'
$ws.Range('E32').Value = '  -1.16%  '
$ws.Range('E33').Value = '  -2.16%  '
$ws.Range('E34').Value = '  -2.53%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.137'
$ws.Range('E35').Value = '  +5.22%  '
$ws.Range('E36').Value = '  +0.15%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.51'
$ws.Range('E37').Value = '  -4.99%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '155.49'
$ws.Range('E38').Value = '  -0.41%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.11'
$ws.Range('E39').Value = '  -2.96%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.366'
$ws.Range('E40').Value = '  -1.99%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.16'
$ws.Range('E41').Value = '  -2.61%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.79'
$ws.Range('E42').Value = '  -2.61%  '
$ws.Range('E43').Value = '  -0.12%  '
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.44'
$ws.Range('E45').Value = '  -4.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.23'
$ws.Range('E46').Value = '  -0.78%  '
$ws.Range('E47').Value = '  -1.80%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.581'
$ws.Range('E48').Value = '  -1.63%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '151.61'
$ws.Range('E49').Value = '  -1.62%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.79'
$ws.Range('E50').Value = '  -2.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.72'
$ws.Range('E51').Value = '  -1.25%  '
